$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- Update "Ativação:" date from 01/01/2016 to 01/01/2023 ---
# (B8/C8 hold the date text, and B13/C13 happen to reuse the exact same shared
#  string in the original workbook, so update all four the same way.)
foreach ($addr in @("B8", "C8", "B13", "C13")) {
    $cell = $ws.Range($addr)
    # Force the value to be stored as text instead of being auto-parsed as a date
    $cell.NumberFormat = "@"
    $cell.Value = "01/01/2023"
}

# Restore the original (non-text) number format / style that the cells had
# before, by pasting formats only from a cell that already carries the
# correct style (column B uses style index 2, column C uses style index 3).
$ws.Range("B10").Copy()
$ws.Range("B8").PasteSpecial($xlPasteFormats)
$ws.Range("B13").PasteSpecial($xlPasteFormats)

$ws.Range("C10").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)
$ws.Range("C13").PasteSpecial($xlPasteFormats)

$ws.Application.CutCopyMode = $false

# --- Add new "Objectives:" body text (row 11) ---
$ws.Range("B10:C10").Copy()
$ws.Range("B11:C11").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

$objectivesText = "Develop knowledge in order to make the student capable of correctly interpreting the technical drawing, knowing the methodologies and tools used in the industry, giving subsidies so that they can execute, interact and modify drawings and projects throughout their professional life."
$ws.Range("B11").Value = $objectivesText
$ws.Range("C11").Value = $objectivesText

# --- Add new "Short syllabus:" body text (row 14) ---
$ws.Range("B13:C13").Copy()
$ws.Range("B14:C14").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

$shortSyllabusText = "Context of the technical drawing in the industry, main tools and techniques used in drawings for the elaboration of projects. Introduction to computer-aided design (CAD)."
$ws.Range("B14").Value = $shortSyllabusText
$ws.Range("C14").Value = $shortSyllabusText

# --- Add new "Syllabus:" body text (row 16) ---
$ws.Range("B15:C15").Copy()
$ws.Range("B16:C16").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

$syllabusText = "Technical drawing standards. Technical terminology and materials for drawing. Perspective representation. Orthogonal design. Scaling and scaling. Cut and section. Auxiliary view and details. Geometric tolerances. Representation of machine elements. Use of software for technical design. Computer-aided design in three dimensions (Solid Modeling). Computer-aided design in two dimensions."
$ws.Range("B16").Value = $syllabusText
$ws.Range("C16").Value = $syllabusText
